# Effect of Charging Stations on EV Market Share
# - "About" sheet: drop the old source citation (page/author/title/URL/year),
#   replace the source reference with "None", remove the old hyperlink, and
#   add a new note explaining the US-specific override (value forced to 0).
# - "EoCSoEVMS" sheet: zero out the "1 more charger per 100k pop" effect size.

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")

# Drop the old citation's hyperlink (used to live on B6: the sciencedirect URL).
foreach ($hl in $about.Hyperlinks) {
    $hl.Delete()
}

# Drop the "last updated" date stamp that used to live in C1.
$about.Range("C1").Clear()

# Remove the old citation block entirely:
#   B4 = 2014 (year)
#   B5 = "The influence of financial incentives ... on electric vehicle adoption"
#   B6 = hyperlink text (sciencedirect URL)
#   B7 = "Page 190, column 2"
# Deleting these rows shifts the old "Notes" block (rows 9-11) up to rows 5-7.
$about.Range("4:7").Delete()

# New note (row 9) explaining the US-specific override. Written before B3 below
# so that new shared strings land in the same order as in the target workbook.
$about.Range("A9").Value = "In the US, we set this to 0 so that increasing EV chargers does not induce additional deployment."

# Replace the old author citation (B3) with "None" (source no longer used).
$about.Range("B3").Value = "None"

$data = $wb.Worksheets.Item("EoCSoEVMS")

# In the US, increasing EV chargers should not induce additional EV share
# deployment, so this effect size is set to 0 (was 0.0012).
$data.Range("B2").Value = 0

# The old citation hyperlink was the only thing using the built-in
# "Hyperlink" cell style; drop the now-unused named style too.
foreach ($cellStyle in $wb.Styles) {
    if ($cellStyle.Name() -eq "Hyperlink") {
        $cellStyle.Delete()
    }
}
